# Update the small lookup table on Sheet1:
#  - swap the two header labels between columns A and B
#  - replace the TiD/basket-item sample data with the new item names / id lists
#  - remove the last (6th) data row that no longer exists in the new data set

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: column A now holds "items", column B now holds "TiD"
$ws.Range("A1").Value = "items"
$ws.Range("B1").Value = "TiD"

# Data rows: column A holds the item name, column B holds the comma separated id list
$ws.Range("A2").Value = "bread"
$ws.Range("B2").Value = "1,2,3,5"

$ws.Range("A3").Value = "sugar"
$ws.Range("B3").Value = "1,2,4,5"

$ws.Range("A4").Value = "milk"
$ws.Range("B4").Value = "1,3,4"

$ws.Range("A5").Value = "powder"
$ws.Range("B5").Value = "1,3,5"

# The sixth row of the old table is no longer part of the data set
$ws.Range("A6:B6").ClearContents()

# Keep the row-height stamp that trails the shrunk table
$ws.Rows(8).RowHeight = 15
